$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear comments that are no longer applicable
$ws.Range("H8").Clear()
$ws.Range("H9").Clear()

# Add new feedstock-related comments
$ws.Range("H15").Value = "corn or cellulosic feedstock only"
$ws.Range("H21").Value = "second generation feedstocks only"
$ws.Range("H24").Value = "grain feedstock only"
$ws.Range("H25").Value = "second generation feedstocks only"

# Update selection
$ws.Range("E23").Select() | Out-Null
